# Update the cached "datetimeFigureOut" date field text from 2025-02-17 to
# 2025-02-26 everywhere it appears: the slide master, every slide layout,
# and the notes master.

$p = $ppt.ActivePresentation
$oldDate = "2025-02-17"
$newDate = "2025-02-26"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout attached to the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Notes master.
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes
